# update scripts wuth new tpm
# Refresh TPM-based NATMI ligand-receptor pair output for Col1a2-Itga11
# (YoungD0 / LR-pairs_lrc2p): recompute existing Sending/Target cluster
# combinations with the new TPM-derived statistics, and add the three
# previously-omitted same-cluster (self-pair) rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col1a2"
$ws.Range("C2").Value = "Itga11"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.042494666666667
$ws.Range("H2").Value = 3.127484
$ws.Range("I2").Value = 0.0007670466909205676
$ws.Range("J2").Value = 0.0007670466909205677
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.006053333333333333
$ws.Range("N2").Value = 0.01816
$ws.Range("O2").Value = 0.0002373030766683641
$ws.Range("P2").Value = 0.0002373030766683641
$ws.Range("Q2").Value = 0.006310567715555555
$ws.Range("R2").Value = 0.05679510944
$ws.Range("S2").Value = 0.0000001820225397037384
$ws.Range("T2").Value = 0.0000001820225397037385
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col1a2"
$ws.Range("C3").Value = "Itga11"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.042494666666667
$ws.Range("H3").Value = 3.127484
$ws.Range("I3").Value = 0.0007670466909205676
$ws.Range("J3").Value = 0.0007670466909205677
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 24.759128
$ws.Range("N3").Value = 74.277384
$ws.Range("O3").Value = 0.970608576546119
$ws.Range("P3").Value = 0.9706085765461191
$ws.Range("Q3").Value = 25.81125889131733
$ws.Range("R3").Value = 232.301330021856
$ws.Range("S3").Value = 0.000744502096818823
$ws.Range("T3").Value = 0.0007445020968188232
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col1a2"
$ws.Range("C4").Value = "Itga11"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.042494666666667
$ws.Range("H4").Value = 3.127484
$ws.Range("I4").Value = 0.0007670466909205676
$ws.Range("J4").Value = 0.0007670466909205677
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7436886666666668
$ws.Range("N4").Value = 2.231066
$ws.Range("O4").Value = 0.02915412037721258
$ws.Range("P4").Value = 0.02915412037721259
$ws.Range("Q4").Value = 0.7752914686604445
$ws.Range("R4").Value = 6.977623217944
$ws.Range("S4").Value = 0.0000223625715620408
$ws.Range("T4").Value = 0.00002236257156204081
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col1a2"
$ws.Range("C5").Value = "Itga11"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1312.703450666667
$ws.Range("H5").Value = 3938.110352
$ws.Range("I5").Value = 0.9658609009611662
$ws.Range("J5").Value = 0.9658609009611662
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.006053333333333333
$ws.Range("N5").Value = 0.01816
$ws.Range("O5").Value = 0.0002373030766683641
$ws.Range("P5").Value = 0.0002373030766683641
$ws.Range("Q5").Value = 7.946231554702222
$ws.Range("R5").Value = 71.51608399231999
$ws.Range("S5").Value = 0.0002292017634317629
$ws.Range("T5").Value = 0.0002292017634317629
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col1a2"
$ws.Range("C6").Value = "Itga11"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1312.703450666667
$ws.Range("H6").Value = 3938.110352
$ws.Range("I6").Value = 0.9658609009611662
$ws.Range("J6").Value = 0.9658609009611662
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.759128
$ws.Range("N6").Value = 74.277384
$ws.Range("O6").Value = 0.970608576546119
$ws.Range("P6").Value = 0.9706085765461191
$ws.Range("Q6").Value = 32501.39276109769
$ws.Range("R6").Value = 292512.5348498792
$ws.Range("S6").Value = 0.9374728742234695
$ws.Range("T6").Value = 0.9374728742234696
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col1a2"
$ws.Range("C7").Value = "Itga11"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1312.703450666667
$ws.Range("H7").Value = 3938.110352
$ws.Range("I7").Value = 0.9658609009611662
$ws.Range("J7").Value = 0.9658609009611662
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7436886666666668
$ws.Range("N7").Value = 2.231066
$ws.Range("O7").Value = 0.02915412037721258
$ws.Range("P7").Value = 0.02915412037721259
$ws.Range("Q7").Value = 976.2426789550259
$ws.Range("R7").Value = 8786.184110595232
$ws.Range("S7").Value = 0.02815882497426484
$ws.Range("T7").Value = 0.02815882497426484
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Col1a2"
$ws.Range("C8").Value = "Itga11"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 45.356022
$ws.Range("H8").Value = 136.068066
$ws.Range("I8").Value = 0.03337205234791334
$ws.Range("J8").Value = 0.03337205234791334
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.006053333333333333
$ws.Range("N8").Value = 0.01816
$ws.Range("O8").Value = 0.0002373030766683641
$ws.Range("P8").Value = 0.0002373030766683641
$ws.Range("Q8").Value = 0.27455511984
$ws.Range("R8").Value = 2.47099607856
$ws.Range("S8").Value = 0.00000791929069689754
$ws.Range("T8").Value = 0.00000791929069689754
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Col1a2"
$ws.Range("C9").Value = "Itga11"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 45.356022
$ws.Range("H9").Value = 136.068066
$ws.Range("I9").Value = 0.03337205234791334
$ws.Range("J9").Value = 0.03337205234791334
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 24.759128
$ws.Range("N9").Value = 74.277384
$ws.Range("O9").Value = 0.970608576546119
$ws.Range("P9").Value = 0.9706085765461191
$ws.Range("Q9").Value = 1122.975554268816
$ws.Range("R9").Value = 10106.77998841934
$ws.Range("S9").Value = 0.03239120022583073
$ws.Range("T9").Value = 0.03239120022583074
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Col1a2"
$ws.Range("C10").Value = "Itga11"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 45.356022
$ws.Range("H10").Value = 136.068066
$ws.Range("I10").Value = 0.03337205234791334
$ws.Range("J10").Value = 0.03337205234791334
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7436886666666668
$ws.Range("N10").Value = 2.231066
$ws.Range("O10").Value = 0.02915412037721258
$ws.Range("P10").Value = 0.02915412037721259
$ws.Range("Q10").Value = 33.730759526484
$ws.Range("R10").Value = 303.576835738356
$ws.Range("S10").Value = 0.0009729328313857052
$ws.Range("T10").Value = 0.0009729328313857053
